$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "anular remito pendiente de fc" (row 56) as finished
$ws.Range("B56").Value = "terminado"

# Add new backlog row
$ws.Range("A58").Value = "FILTRAR REMITO PENDIENTE DE FC X CLIENTE"
$ws.Range("B58").Value = "no comenzado"

# Update the active selection to match the new last row
$ws.Range("A57").Select()
